$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.401.92'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '1.549.27'
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.76%  '
$ws.Range("E6").Value = '  -2.24%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.01'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.67%  '
$ws.Range("E9").Value = '  -1.94%  '
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("D12").Value = '1.771.88'
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("D13").Value = '1.550.88'
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").Value = '28.385.31'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("E15").Value = '  -2.80%  '
$ws.Range("E16").Value = '  -2.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.33'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.13%  '
$ws.Range("D20").Value = '0.0₃0672'
$ws.Range("E20").Value = '  -2.99%  '
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.16%  '
$ws.Range("E24").Value = '  -3.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("E26").Value = '  -2.16%  '
$ws.Range("E27").Value = '  -1.58%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.54%  '
$ws.Range("E30").Value = '  -3.57%  '
$ws.Range("E31").Value = '  -5.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.64%  '
$ws.Range("D33").Value = '1.381.37'
$ws.Range("E33").Value = '  -1.34%  '
$ws.Range("E34").Value = '  -3.85%  '
$ws.Range("E35").Value = '  +0.81%  '
$ws.Range("E36").Value = '  -3.79%  '
$ws.Range("E37").Value = '  -1.93%  '
$ws.Range("E38").Value = '  -3.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0161'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.57%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("B42").Value = 'ImmutableX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.508'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.768'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0455'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.08%  '
$ws.Range("E45").Value = '  -2.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.22%  '
$ws.Range("D47").Value = '1.684.21'
$ws.Range("E47").Value = '  -2.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.872'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.36%  '
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '43.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.09%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '85.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("D51").Value = '0.0₆0103'
$ws.Range("E51").Value = '  -0.65%  '
